$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old A1:D60 data (layout changes from 4 columns to 2, and one row is inserted)
$ws.Range("A1:D60").ClearContents()

$ws.Range("A1").Value = "アーコ"
$ws.Range("B1").Value = "Ako"

$ws.Range("A2").Value = "MP_SET_MOVIE e8_Ta"
$ws.Range("B2").Value = "MP_SET_MOVIE e8_Ta"

$ws.Range("A3").Value = "MP_SET_LOOP 6 on"
$ws.Range("B3").Value = "MP_SET_LOOP 6 on"

$ws.Range("A4").Value = "\n<アーコ>やぁやぁ。アーコだよ。`n何か欲しい物とかあるの？`n見てく？"
$ws.Range("B4").Value = "\n<Ako>Salutations. Ako here!`nAre you here to buy something?`nWant to look?"

$ws.Range("A5").Value = "\n<アーコ>武器になりそうなものとか見つけた？`nまずはそういうの探すといいかもね。"
$ws.Range("B5").Value = "\n<Ako>Have you found anything you could use as a weapon?`nIf not, you should probably find something like that first."

$ws.Range("A6").Value = "\n<アーコ>何食べたらそんなにおっぱい大きくなるの？`n横縞模様だからおっぱい大きく見えるの？`n錯覚？"
$ws.Range("B6").Value = "\n<Ako>What do you eat to make your boobs so big?`nIs it the stripes that make them look bigger?`nIs it some kind of illusion?"

$ws.Range("A7").Value = "\n<アーコ>囚人服に合ってるよ。`nまさに囚人って感じ。`n褒めてないよ。"
$ws.Range("B7").Value = "\n<Ako>That prison outfit suits you.`nI mean, you are a prisoner I guess.`nDon't take it as a compliment."

$ws.Range("A8").Value = "\n<アーコ>あたしは淫魔だし、気持ちが分かるから`nどっちの味方もしないよ。`nでも仲良くしたいって思ってる。"
$ws.Range("B8").Value = "\n<Ako>I'm a succubus, and I know how you feel.`nI won't pick sides.`nBut I hope we can be friends."

$ws.Range("A9").Value = "\n<アーコ>シィナは猫の時からの友達だよ。`nリリーに淫魔にされたんだってさ。`nあたしは違うやつに淫魔にされたけど。"
$ws.Range("B9").Value = "\n<Ako>Shina has been my friend since we were just cats!`nLily turned her into a succubus, and now she's like that.`nA different person turned me into a succubus though."

$ws.Range("A10").Value = "\n<アーコ>あたしも淫魔だから精液が一番美味しいと思うけど`nカリカリもかつおぶしも同じくらい美味しいと思う。`n一番がいっぱいあるといいね。"
$ws.Range("B10").Value = "\n<Ako>I'm a succubus as well so of course semen is the most delicious`nthing in the world... But I love crisps and katsuobushi as well.`nI can never get enough of them~"

$ws.Range("A11").Value = "\n<アーコ>ここは元々大金持ちが住んでたんだって。`nメイドも使用人もいっぱい居たってさ。`nなんでこんな森の中に？って感じ。"
$ws.Range("B11").Value = "\n<Ako>A very rich person used to live here...`nAnd lots of maids and servants too. It makes me think-`nWhat was it doing out in the forest?"

$ws.Range("A12").Value = "\n<アーコ>みんなお金に興味ないから`nあたしがこの館から色々持って行っても何とも思わないんだよ。`nだから、ここで仕入れて町で売ったりする。"
$ws.Range("B12").Value = "\n<Ako>Because no one has any interest in gold here,`nthey don't care if I take stuff from the mansion.`nSo, I just take it to town and sell it."

$ws.Range("A13").Value = "\n<アーコ>その見返りにタバコとか持ってきてあげるの。`nまぁ、それもお金取るけどね。"
$ws.Range("B13").Value = "\n<Ako>I'll give you cigars and stuff for something good.`n...Well, for the right price that is."

$ws.Range("A14").Value = "\n<アーコ>淫魔の巣に閉じ込められるのは`n人間にとって辛いでしょ。`n少しでもえっちな気分になるとすぐ硬くなっちゃう。"
$ws.Range("B14").Value = "\n<Ako>It must be so difficult being a human `nlocked up in the succubus hideout.`nEven being just a bit aroused will make you hard."

$ws.Range("A15").Value = "\n<アーコ>行き詰ったらとりあえず攻撃してみたらいいよ。`n壊せるものとかあるかもしれないし。"
$ws.Range("B15").Value = "\n<Ako>If things look hopeless, just keep attacking.`nYou might just make a break through."

$ws.Range("A16").Value = "\n<アーコ>お兄ちゃんも下に居るの？`n面白そうだから後で行ってみよー。"
$ws.Range("B16").Value = "\n<Ako>Your big brother is down there too right?`nHow interesting, I'll be sure to give him a visit later-."

$ws.Range("A17").Value = "\n<アーコ>妹も捕まってるの？`n後で行ってみよー。`n仲良くなれるかな？"
$ws.Range("B17").Value = "\n<Ako>So your sister has been captured as well?`nI'll go and visit her later-.`nI wonder if we would get along well?"

$ws.Range("A18").Value = "買い物"
$ws.Range("B18").Value = "I'd like to buy something."

$ws.Range("A19").Value = "用事はない"
$ws.Range("B19").Value = "I don't need anything right now."

$ws.Range("A20").Value = "\n<アーコ>何を買ってくれるのかな？"
$ws.Range("B20").Value = "\n<Ako>What would you like to buy?"

$ws.Range("A21").Value = "\n<アーコ>ははーん。さては冷やかしだな？`nふーん！"
$ws.Range("B21").Value = "\n<Ako>Yeah yeah. Did you come just to stare?`nHmmph!"

$ws.Range("A22").Value = "\n<アーコ>まいどあり～！"
$ws.Range("B22").Value = "\n<Ako>Thank you for stopping by~!"

$ws.Range("A23").Value = "\n<アーコ>ふーん。"
$ws.Range("B23").Value = "\n<Ako>Hmmph."

$ws.Range("A24").Value = "リリー"
$ws.Range("B24").Value = "Lily"

$ws.Range("A25").Value = "シィナ"
$ws.Range("B25").Value = "Shina"

$ws.Range("A26").Value = "アーコイベ用"
$ws.Range("B26").Value = "アーコイベ用"

$ws.Range("A27").Value = "ライム"
$ws.Range("B27").Value = "Lime"

$ws.Range("A28").Value = "自動実行カギ"
$ws.Range("B28").Value = "自動実行カギ"

$ws.Range("A29").Value = "\n<リリー>あ、シィナあんた。`nちゃんと鍵掛けた？"
$ws.Range("B29").Value = "\n<Lily>Oh Shina, you're back.`nDid you lock them up?`n"

$ws.Range("A30").Value = "\n<シィナ>は？知らんし。"
$ws.Range("B30").Value = "\n<Shina>Huh, me?`nBeats me."

$ws.Range("A31").Value = "\n<リリー>ちょっと！`nしっかりしてよね！`n逃げられたらどうすんのよ！"
$ws.Range("B31").Value = "\n<Lily>Whoa whoa whoa!`nYou need to take this seriously!`nWhat if they got out!"

$ws.Range("A32").Value = "\n<リリー>はぁ～～～・・・`nほんとあんたっていっつも・・・"
$ws.Range("B32").Value = "\n<Lily>Sigh....`nWhy are you always so... so..."

$ws.Range("A33").Value = "\n<シィナ>はぁ？`n知らんにゃ。`n鍵持ってないし。アタシじゃねー！"
$ws.Range("B33").Value = "\n<Shina>Whaa?`nI dunno-nya.`nI don't have the key. It's not my fault!"

$ws.Range("A34").Value = "\n<ライム>私も持ってないよ。"
$ws.Range("B34").Value = "\n<Lime>I don't have a key either."

$ws.Range("A35").Value = "\n<リリー>ん？"
$ws.Range("B35").Value = "\n<Lily>Hm?"

$ws.Range("A36").Value = "\n<リリー>・・・！！"
$ws.Range("B36").Value = "\n<Lily>...!!"

$ws.Range("A37").Value = "\n<リリー>私が持ってた。`n・・・`n鍵、閉めたっけ？"
$ws.Range("B37").Value = "\n<Lily>I have the key.`n...`nDid you at least shut the door?"

$ws.Range("A38").Value = "\n<シィナ>知らんし。`nアホにゃ。"
$ws.Range("B38").Value = "\n<Shina>I dunno.`nWhy don't you do it yourself-nya."

$ws.Range("A39").Value = "\n<ライム>とりあえず、確認しに戻ろっか。`nリリー。"
$ws.Range("B39").Value = "\n<Lime>Well for now, let's go back and check,`nLily?"

$ws.Range("A40").Value = "\n<シィナ>ごめんなさいは～？"
$ws.Range("B40").Value = "\n<Shina>Now who should be the one apologizing?~"

$ws.Range("A41").Value = "\n<リリー>えへっ♥"
$ws.Range("B41").Value = "\n<Lily>Eei♥"

$ws.Range("A42").Value = "最初のイベント"
$ws.Range("B42").Value = "最初のイベント"

$ws.Range("A43").Value = "\n<ライム>やっぱりね！"
$ws.Range("B43").Value = "\n<Lime>Just as I thought！"

$ws.Range("A44").Value = "\n<ライム>逃げる時間そんなにないよなー、って思ってたの。`nどこかに隠れてたのかな？`n待ってたら来ると思ったよー！"
$ws.Range("B44").Value = "\n<Lime>They couldn't have had much time to escape, I thinks.`nMaybe they're hiding somewhere?`nIf we wait, I'm sure they'll come out!~"

$ws.Range("A45").Value = "ーーーーー基本変更点ーーーーー"

$ws.Range("A46").Value = "\n<ライム>にへへへへー♥`nつーかまーえた♥"
$ws.Range("B46").Value = "\n<Lime>Ehehehehe-♥`nCau~ght you~♥"

$ws.Range("A47").Value = "\n<\n[3]>ぬるぬるして動きにくいでしょー。`n早く逃げないと白いの出させちゃうぞー？"
$ws.Range("B47").Value = "\n<\n[3]>It's hard to move when you're so slimy, right?`nIf you don't get away quickly, you're gonna end up`nleaking out white stuff you know?"

$ws.Range("A48").Value = "\C[3]※捕まるとタイミングバーが表示されます。`n\C[0]タイミングよく黄か赤で止めてください。`n赤で止めると被ダメージが半減します。"
$ws.Range("B48").Value = "\C[3]※A timing bar will show if you are caught.`n\C[0]Time your key press to match either the yellow or red areas.`nIf you stop in the red area, you will take half damage."

$ws.Range("A49").Value = "\n<\n[3]>むにゅー・・・♥`n柔らかくて溶けちゃいそうでしょー♥`n気持ちいい気持ちいいー♥"
$ws.Range("B49").Value = "\n<\n[3]>Squeeze---...♥`nThey feel so soft you could melt right?♥`nSo plush, so pleasureable-♥"

$ws.Range("A50").Value = "\n<\n[3]>あれー？もう出ちゃうのー？`nおっぱい我慢できなかったー？あはは♥`nじゃあ一回だけ、出しちゃおっかー♪"
$ws.Range("B50").Value = "\n<\n[3]>Oh-? Are you gonna squirt already?`nCouldn't you resist my boobies? Aha!♥`nAlright then, I'm sure you can squirt just a little-♪"

$ws.Range("A51").Value = "\n<\n[3]>あっあっ♥おっぱいの間でぴくぴくしてるー♥`nおちんちん喜んでくれたみたい♥`nうれしー♥"
$ws.Range("B51").Value = "\n<\n[3]>Ahaha♥ I can feel it twitching in my cleavage♥`nYour penis seems really happy♥`nI'm happy too♥"

$ws.Range("A52").Value = "\C[1]SAN値が1下がった・・・（現在SAN値\v[270]）"
$ws.Range("B52").Value = "\C[1]Sanity decreased by 1... (Current Sanity: \v[270]）"

$ws.Range("A53").Value = "\n<ライム>もうー。ひょっとして全然抵抗する気ないのー？`nゲームはまだ始まったばかりなのに。`n次は本気で搾っちゃうからねー？"
$ws.Range("B53").Value = "\n<Lime>Gee~ Were you trying to resist at all?`nAnd the game only just started too.`nNext time, I'll milk you for real, okay?"

$ws.Range("A54").Value = "\n<ライム>あっ！逃げられちゃった・・・`nまぁいいっかー。"
$ws.Range("B54").Value = "\n<Lime>Ack! They got away...`nOh well..."

$ws.Range("A55").Value = "\n<ライム>えーっと、リリーの行ってたこと、`n聞こえてたかな？`n脱出ゲームがどーのこーの・・・"
$ws.Range("B55").Value = "\n<Lime>So, um... Did you hear about what`nLily was going on about?`nAbout this...escape game and such..."

$ws.Range("A56").Value = "\n<ライム>この館は広いからねー。`nなかなか出られないと思うけど・・・`n頑張って色んなところ、探してみてね。"
$ws.Range("B56").Value = "\n<Lime>This mansion is really quite huge you know~`nI don't think you'll be able to leave easily...`nBut do everything you can, and be sure to search around OK?"

$ws.Range("A57").Value = "\n<ライム>リリーもシィナも楽しそうだから。`n簡単に諦めたりしちゃーダメだよー？`nいっぱい遊ばれてね。"
$ws.Range("B57").Value = "\n<Lime>Lili and Shina seem to like the idea.`nSo, you'd better not give up so easily, alright?`nLet's play lots and lots OK?"

$ws.Range("A58").Value = "\n<ライム>次は一回射精したぐらいじゃ`n許してあげないからねー？くすくす♥`nばいばーい♥"
$ws.Range("B58").Value = "\n<Lime>Next time, I won't let you get away with`njust one ejaculation, okay? Teehee♥`nBye bye♥"

$ws.Range("A59").Value = "\n<ライム>じゃー頑張ってねー。`n次は私も本気でぴゅっぴゅさせにいっちゃうからねー。"
$ws.Range("B59").Value = "\n<Lime>Now then- Do your best~`nNext time I'll make you go `"pew pew`" for realsies~"

$ws.Range("A60").Value = "シィナイベ用"
$ws.Range("B60").Value = "Shinaイベ用"

$ws.Range("A61").Value = "梯子"
$ws.Range("B61").Value = "梯子"

# Reset any auto-applied custom row heights so rows stay at default height,
# matching original formatting (content changes only, no row-height changes).
$ws.Range("A1:A61").EntireRow.AutoFit()
